# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values change for rows 2-23; update each cell with its new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 3
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 2
    16 = 2
    17 = 0
    18 = 0
    19 = 1
    20 = 3
    21 = 4
    22 = 1
    23 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
